# "se modifico controlador excel"
#
# The sheet had a stray extra column (W2 = 9) beyond the real data table
# (which ends at column V, FECHA_EMI_CERTF). The fix:
#   - V2 becomes the text "-" instead of the number 8.
#   - W2 (and the whole now-unused column W) is removed, so the sheet's
#     used range shrinks from A1:W2 to A1:V2.
#   - The view is scrolled right to column H and the selection left on V3,
#     matching where the editor's cursor ended up after trimming the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# V2: number 8 -> text "-" (shared string).
$ws.Range("V2").Value = "-"

# Drop column W entirely (it only held the orphan value 9 in W2), which
# also updates the sheet dimension/row spans to A1:V2 / 1:22.
$ws.Range("W:W").Delete() | Out-Null

# Best-effort view state matching: scroll so column H is left-most visible
# and leave the selection on V3.
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("V3").Select() | Out-Null
